$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 values (columns B..AH) to 2 decimal places (custom accuracy).
# These are literal rounded values (matching the target dataset), not a
# live ROUND() formula, so we assign them directly.
$row5 = @{
    "B5"  = 17.78
    "C5"  = 13.02
    "D5"  = 1.15
    "E5"  = 38.36
    "F5"  = 31.74
    "G5"  = 13.99
    "H5"  = 54.6
    "I5"  = 21.52
    "J5"  = 9.51
    "K5"  = 14.25
    "L5"  = 15.47
    "M5"  = 16.2
    "N5"  = 4.47
    "O5"  = 13.91
    "P5"  = 19.77
    "Q5"  = 11.75
    "R5"  = 0.85
    "S5"  = 0.75
    "T5"  = 204.41
    "U5"  = 38.9
    "V5"  = 12.84
    "W5"  = 26.1
    "X5"  = 13.89
    "Y5"  = 1.73
    "Z5"  = 26.32
    "AA5" = 11.34
    "AB5" = 10.1
    "AC5" = 11.86
    "AD5" = 16.2
    "AE5" = 0.5600000000000001
    "AF5" = 49.34
    "AG5" = 7.21
    "AH5" = 16.05
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# Remove row 6 entirely (dataset trimmed from 6 to 5 rows).
$ws.Rows.Item(6).Delete()
